# Generate Report for Handoff
#
# The localization status report moves from "In Translation" to
# "Ready for handoff": the Status cells on the Overview, zh-cn and de-de
# sheets are updated, the corresponding "generate"/"handoff" timestamps are
# bumped to the new report time, and the Status columns are widened to fit
# the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bumped to the new handoff generation time ---
$wsOverview.Range("G2").Value = "2016-08-13 14:47:26"
$wsDeDe.Range("H2").Value = "2016-08-13 14:47:26"
$wsZhCn.Range("H2").Value = "2016-08-13 14:47:19"

# --- Widen the Status columns to fit "Ready for handoff" ---
$wsOverview.Range("E1:F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
